# Update the report title strings in row 1 from "10.b.1" to "10.b.1.1"
# (Kyrgyz title in A1, English title in C1; the Russian title in B1 is
# already correct and untouched), then restore the active-cell selection
# that Excel recorded on last save (L8).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "10.b.1.1 Агымдардын түрлөрү жана алуучу өлкөлөр жана донор-өлкөлөр боюнча бөлунүшүндөгү  өнүктүрүү максатында ресурстар агымынын жалпы көлөмү"
$ws.Range("C1").Value = "10.b.1.1 Total resource flows for development, by recipient and donor countries and type of flow (e.g. official development assistance, foreign direct investment and other flows)"

$ws.Range("L8").Select()
